$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.336.96'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.933.93'
$ws.Range("E3").Value = '  +0.19%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7553'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.87%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.86'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.74%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("B8").Value = 'Cardano'
$ws.Range("C8").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3171'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.70%  '
$ws.Range("B9").Value = 'Solana'
$ws.Range("C9").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '27.54'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06994'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.66%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7769'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.58%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08017'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.84%  '
$ws.Range("D13").Value = '1.933.83'
$ws.Range("E13").Value = '  +0.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.342'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.57%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '94.39'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.37'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.91%  '
$ws.Range("D17").Value = '30.349.34'
$ws.Range("E17").Value = '  +0.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '253.31'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.93%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007926'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.28%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.729'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.89%  '
$ws.Range("D21").Value = '2.191.10'
$ws.Range("E21").Value = '  +0.59%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.0000'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.09%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.19%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.653'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.71%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.466'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.90'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("B27").Value = 'Stellar'
$ws.Range("C27").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1335'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +4.00%  '
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.99'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.09%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.197'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -5.01%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.366'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.42%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.516'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.75%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.392'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.73%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.121'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05144'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.28%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.281'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.54%  '
$ws.Range("E36").Value = '  +0.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.771'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.36%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01957'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.14%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.799'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '77.31'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.417'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.04%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4452'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.46%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.961'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.001'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8331'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.56'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.19%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.760'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.472'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.95%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.34%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '978.53'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1180'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.09%  '
